$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows before row 429, shifting existing rows 429-450
# down to 431-452 (this also naturally moves old rows 449-450 to 451-452,
# matching the appended rows at the bottom of the table).
$ws.Range("A429:A430").EntireRow.Insert()

# Fill in the new row 429 (Primera) with the latest weekly data point.
$ws.Cells.Item(429, 1).Value = 11
$ws.Cells.Item(429, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(429, 3).Value = "Bíobío"
$ws.Cells.Item(429, 4).Value = 45147
$ws.Cells.Item(429, 5).Value = 8
$ws.Cells.Item(429, 6).Value = 100112009
$ws.Cells.Item(429, 7).Value = "Acelga"
$ws.Cells.Item(429, 8).Value = "Sin especificar"
$ws.Cells.Item(429, 9).Value = "Primera"
$ws.Cells.Item(429, 10).Value = 200
$ws.Cells.Item(429, 11).Value = 600
$ws.Cells.Item(429, 12).Value = 700
$ws.Cells.Item(429, 13).Value = 650
$ws.Cells.Item(429, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(429, 15).Value = "Región de Ñuble"
$ws.Cells.Item(429, 16).Value = 650
$ws.Cells.Item(429, 17).Value = 1
$ws.Cells.Item(429, 18).Value = "Hortaliza"
$ws.Cells.Item(429, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Fill in the new row 430 (Segunda) with the latest weekly data point.
$ws.Cells.Item(430, 1).Value = 11
$ws.Cells.Item(430, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(430, 3).Value = "Bíobío"
$ws.Cells.Item(430, 4).Value = 45147
$ws.Cells.Item(430, 5).Value = 8
$ws.Cells.Item(430, 6).Value = 100112009
$ws.Cells.Item(430, 7).Value = "Acelga"
$ws.Cells.Item(430, 8).Value = "Sin especificar"
$ws.Cells.Item(430, 9).Value = "Segunda"
$ws.Cells.Item(430, 10).Value = 100
$ws.Cells.Item(430, 11).Value = 500
$ws.Cells.Item(430, 12).Value = 500
$ws.Cells.Item(430, 13).Value = 500
$ws.Cells.Item(430, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(430, 15).Value = "Región de Ñuble"
$ws.Cells.Item(430, 16).Value = 500
$ws.Cells.Item(430, 17).Value = 1
$ws.Cells.Item(430, 18).Value = "Hortaliza"
$ws.Cells.Item(430, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
